$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Materials"
$ws.Range("C1").Value = "Industrials"
$ws.Range("D1").Value = "Consumer Discretionary"
$ws.Range("E1").Value = "Consumer Staples"
$ws.Range("F1").Value = "Health Care"
$ws.Range("G1").Value = "Financials"
$ws.Range("H1").Value = "Information Technology"
$ws.Range("I1").Value = "Telecommunication Services"
$ws.Range("J1").Value = "Utilities"
$ws.Range("K1").Value = "Real Estate"

$ws.Range("A2").Value = 42614
$ws.Range("A3").Value = 42705
$ws.Range("A4").Value = 42795
$ws.Range("A5").Value = 42887

$ws.Range("A2:A5").NumberFormat = "m/d/yyyy"
$ws.Range("B2:K5").Style = "Percent"
$ws.Range("B2:K5").NumberFormat = "0.00%"
